# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook

# 1) Rename the "Include from RoleClass" sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from RoleClass")
$wsInclude.Name = "Include #0"

# 2) Update the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8)
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), pushing
# Description / Purpose / Copyright / Immutable down by one row (11-14 -> 12-15).

# Capture the existing rows 11-14 values before they get overwritten.
$oldRows = @()
for ($r = 11; $r -le 14; $r++) {
    $oldRows += ,@($wsMeta.Cells.Item($r, 1).Value(), $wsMeta.Cells.Item($r, 2).Value())
}

# Create the new row 15 by copying the format (style) of row 14, so the
# newly-extended row matches the sheet's normal-data-row style instead of
# picking up a blank/default style.
$wsMeta.Range("A14:B14").Copy()
$wsMeta.Range("A15:B15").PasteSpecial(-4122)

# Write the captured values back out, shifted down by one row.
for ($i = 0; $i -lt $oldRows.Length; $i++) {
    $targetRow = 12 + $i
    $wsMeta.Cells.Item($targetRow, 1).Value = $oldRows[$i][0]
    $wsMeta.Cells.Item($targetRow, 2).Value = $oldRows[$i][1]
}

# Finally, write the new Jurisdiction row into the now-vacated row 11.
$wsMeta.Cells.Item(11, 1).Value = "Jurisdiction"
$wsMeta.Cells.Item(11, 2).Value = ""
